$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1156.6111
$ws.Range("I33").Value = 238.6875
$ws.Range("K33").Value = 238.6875
$ws.Range("M33").Value = -9.6875

$ws.Range("H40").Value = 1436.6666
$ws.Range("I40").Value = 1397
$ws.Range("J40").Value = 1471.375
$ws.Range("K40").Value = 1397
$ws.Range("L40").Value = 1471.375
$ws.Range("M40").Value = -1222
$ws.Range("N40").Value = -1821.375

$ws.Range("H70").Value = 1366.5555
$ws.Range("I70").Value = 1412.375
$ws.Range("K70").Value = 4237.125
$ws.Range("M70").Value = -3967.125

$ws.Range("H73").Value = 1366.5555
$ws.Range("I73").Value = 1412.375
$ws.Range("K73").Value = 4237.125
$ws.Range("M73").Value = -3301.125

$ws.Range("H86").Value = 68992.8
$ws.Range("I86").Value = 113398.78
$ws.Range("J86").Value = 2383.8333
$ws.Range("K86").Value = 113398.78
$ws.Range("L86").Value = 2383.8333
$ws.Range("M86").Value = -112275.78
$ws.Range("N86").Value = -4629.8333

$ws.Range("H89").Value = 68992.8
$ws.Range("I89").Value = 113398.78
$ws.Range("J89").Value = 2383.8333
$ws.Range("K89").Value = 566993.9
$ws.Range("L89").Value = 11919.1665
$ws.Range("M89").Value = -561377.9
$ws.Range("N89").Value = -23151.1665

$ws.Range("H104").Value = 20000132
$ws.Range("I104").Value = 20000132
$ws.Range("K104").Value = 60000396
$ws.Range("M104").Value = -59998649

$ws.Range("H111").Value = 5271148
$ws.Range("J111").Value = 20002198
$ws.Range("L111").Value = 60006594
$ws.Range("N111").Value = -60012728

$ws.Range("H116").Value = 3197.6667
$ws.Range("J116").Value = 3197.6667
$ws.Range("L116").Value = 3197.6667
$ws.Range("N116").Value = -10081.6667

$ws.Range("H132").Value = 4102436.5
$ws.Range("I132").Value = 4549718.5
$ws.Range("J132").Value = 2350
$ws.Range("K132").Value = 13649155.5
$ws.Range("L132").Value = 7050
$ws.Range("M132").Value = -13646625.5
$ws.Range("N132").Value = -12110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1770.359
$ws.Range("I132").Value = 1786.2059
$ws.Range("K132").Value = 5358.6177
$ws.Range("M132").Value = -2828.6177

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 59172.61
$ws.Range("I20").Value = 75550.5
$ws.Range("J20").Value = 1850
$ws.Range("K20").Value = 75550.5
$ws.Range("L20").Value = 1850
$ws.Range("M20").Value = -75303.5
$ws.Range("N20").Value = -2344

$ws.Range("H86").Value = 114088.6
$ws.Range("I86").Value = 141813.25
$ws.Range("J86").Value = 3190
$ws.Range("K86").Value = 141813.25
$ws.Range("L86").Value = 3190
$ws.Range("M86").Value = -140690.25
$ws.Range("N86").Value = -5436

$ws.Range("H89").Value = 114088.6
$ws.Range("I89").Value = 141813.25
$ws.Range("J89").Value = 3190
$ws.Range("K89").Value = 709066.25
$ws.Range("L89").Value = 15950
$ws.Range("M89").Value = -703450.25
$ws.Range("N89").Value = -27182

$ws.Range("H105").Value = 126882.125
$ws.Range("I105").Value = 201855.8
$ws.Range("J105").Value = 1926
$ws.Range("K105").Value = 201855.8
$ws.Range("L105").Value = 1926
$ws.Range("M105").Value = -200108.8
$ws.Range("N105").Value = -5420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 48061
$ws.Range("J57").Value = 48061
$ws.Range("L57").Value = 48061
$ws.Range("N57").Value = -49181

$ws.Range("H107").Value = 538.78125
$ws.Range("I107").Value = 505.375
$ws.Range("J107").Value = 639
$ws.Range("K107").Value = 505.375
$ws.Range("L107").Value = 639
$ws.Range("M107").Value = 1414.625
$ws.Range("N107").Value = -4479

$ws.Range("H132").Value = 3182.2273
$ws.Range("I132").Value = 2971.0303
$ws.Range("J132").Value = 3815.818
$ws.Range("K132").Value = 8913.090899999999
$ws.Range("L132").Value = 11447.454
$ws.Range("M132").Value = -6383.090899999999
$ws.Range("N132").Value = -16507.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 268.5
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = 30
$ws.Range("M15").Value = 110

$ws.Range("H131").Value = 800.27
$ws.Range("J131").Value = 800.27
$ws.Range("L131").Value = 2400.81
$ws.Range("N131").Value = -12480.81

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1289.4445
$ws.Range("I122").Value = 951.25
$ws.Range("J122").Value = 1560
$ws.Range("K122").Value = 2853.75
$ws.Range("L122").Value = 4680
$ws.Range("M122").Value = -403.75
$ws.Range("N122").Value = -9580

$ws.Range("H132").Value = 2751.75
$ws.Range("I132").Value = 2628.162
$ws.Range("J132").Value = 4276
$ws.Range("K132").Value = 7884.485999999999
$ws.Range("L132").Value = 12828
$ws.Range("M132").Value = -5354.485999999999
$ws.Range("N132").Value = -17888

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 12043.75
$ws.Range("J42").Value = 12043.75
$ws.Range("L42").Value = 12043.75
$ws.Range("N42").Value = -13169.75

$ws.Range("H49").Value = 12043.75
$ws.Range("J49").Value = 12043.75
$ws.Range("L49").Value = 12043.75
$ws.Range("N49").Value = -12337.75

$ws.Range("H68").Value = 4265.5
$ws.Range("I68").Value = 2350.5
$ws.Range("J68").Value = 4812.643
$ws.Range("K68").Value = 2350.5
$ws.Range("L68").Value = 4812.643
$ws.Range("M68").Value = -1601.5
$ws.Range("N68").Value = -6310.643

$ws.Range("H71").Value = 4265.5
$ws.Range("I71").Value = 2350.5
$ws.Range("J71").Value = 4812.643
$ws.Range("K71").Value = 11752.5
$ws.Range("L71").Value = 24063.215
$ws.Range("M71").Value = -8008.5
$ws.Range("N71").Value = -31551.215

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 6981
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 6981
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 6981
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -7441

$ws.Range("H132").Value = 2733.9714
$ws.Range("I132").Value = 2936.5173
$ws.Range("J132").Value = 1755
$ws.Range("K132").Value = 8809.5519
$ws.Range("L132").Value = 5265
$ws.Range("M132").Value = -6279.5519
$ws.Range("N132").Value = -10325
